$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = "District Court"
$ws.Range("C2").Value = "Original Suit(OS)"
$ws.Range("D2").Value = "District Court"
